$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nameMap = @{}
$nameMap["Michael Walton"] = "Lori Lowery"
$nameMap["Samuel Lester"] = "Pamela Clayton"
$nameMap["Jeffrey Campos"] = "Darren Anderson"
$nameMap["Christopher Garrett"] = "Samantha Spencer"
$nameMap["Melissa Welch"] = "Keith Chavez"
$nameMap["Jacob Hutchinson"] = "Timothy Phillips"
$nameMap["Alicia Graham"] = "Jonathan Aguilar"
$nameMap["Kaitlin Williamson"] = "Kristen Riley"
$nameMap["Andrew Roman"] = "Justin Henderson"
$nameMap["Scott Orozco"] = "Michael Mendez"
$nameMap["Joseph Lee"] = "Christopher Park"
$nameMap["Eileen Hill"] = "Anthony Gill"
$nameMap["Victor Fitzgerald"] = "William Gonzalez"
$nameMap["Raymond Fuller"] = "Tiffany Wade"
$nameMap["April Flores"] = "Patricia Miller"
$nameMap["Robert Wright"] = "Jennifer Roach"
$nameMap["James Becker"] = "Billy Brewer"
$nameMap["Travis Villarreal"] = "Victoria Silva"
$nameMap["Dylan Schultz"] = "Denise Smith"
$nameMap["Adam Bennett"] = "Carl Edwards"
$nameMap["Erica Thomas"] = "Charles Terrell"
$nameMap["Eric Cobb"] = "Heather Woods"
$nameMap["Michael Briggs"] = "Alexis Wright"
$nameMap["Richard Clark"] = "Robert Young"
$nameMap["Michael Taylor"] = "Tracey Solomon"
$nameMap["Jill Travis"] = "Shawn Ray"
$nameMap["Breanna Jensen"] = "Sarah Johnson"
$nameMap["Karen Bowman"] = "Jeremy Peterson"
$nameMap["Joshua Cox"] = "Monique Reid"
$nameMap["Sandra Hobbs"] = "Anne Wilkerson"
$nameMap["Kevin Lee"] = "Amanda Peters"
$nameMap["Deanna Blair"] = "Michelle Lynch"
$nameMap["Michelle Williams"] = "Abigail Orozco"
$nameMap["Anthony Jones"] = "Ian Anderson"
$nameMap["David Jensen"] = "Anthony Martin"
$nameMap["Kenneth Mullen"] = "Frank Mcbride"
$nameMap["Matthew Robinson"] = "Heather Smith"
$nameMap["Ashley Hopkins"] = "William Becker"
$nameMap["Aaron Nichols"] = "Robert Sanchez"
$nameMap["Lauren Tyler"] = "Mark Brown"
$nameMap["Jeffrey Carlson"] = "Shelley Wright"
$nameMap["Andre Howard"] = "Beth Morgan"
$nameMap["Stephen Suarez"] = "Joseph Ray"
$nameMap["Andrew Sullivan"] = "David Suarez"
$nameMap["Maria Meza"] = "Tina Reeves"
$nameMap["Roberta Jenkins"] = "Erika Williams"
$nameMap["Meghan Dunn"] = "William Oneill"
$nameMap["Desiree Brock"] = "Scott Michael"
$nameMap["Laura Watson"] = "Paula Smith"
$nameMap["Valerie Cohen"] = "Raymond Park"
$nameMap["Kimberly Christensen"] = "Robert Calderon"
$nameMap["Brian Jones"] = "Michael Flowers"
$nameMap["Monica Olsen"] = "Yvonne Bennett"
$nameMap["Caitlin Flores"] = "Kim Hale"
$nameMap["Wendy Waters"] = "Lisa Warren"
$nameMap["Kristina Torres"] = "Victoria Young"
$nameMap["Crystal Nielsen"] = "Tiffany Jones"
$nameMap["Sarah Jordan"] = "Jessica Carney"
$nameMap["Austin Solomon"] = "Charles Thompson"
$nameMap["Juan Rodriguez"] = "Cameron Patterson"
$nameMap["Jennifer Gibson"] = "Amy Matthews"
$nameMap["Linda Miller"] = "Christopher Henderson"
$nameMap["Victor Martinez"] = "Jake Henderson"
$nameMap["Andrea Bryan"] = "Joanne Davis"
$nameMap["Kerry Day"] = "Mary Nichols"
$nameMap["Jesse Flores"] = "Joseph Jensen"
$nameMap["Michael Grant"] = "Blake Park"
$nameMap["Michael Simmons"] = "Kevin Parrish"
$nameMap["Scott Morales"] = "Christian Brown"
$nameMap["Jose Medina"] = "James Lane"
$nameMap["Nina Murphy"] = "Jennifer Waller"
$nameMap["Brittany Nelson"] = "Jennifer Kline"
$nameMap["Dr. Patricia Gill"] = "Christine Patrick"
$nameMap["Sabrina Patterson"] = "Morgan Martinez"
$nameMap["Morgan Hernandez"] = "Theresa Black"
$nameMap["Whitney Powell"] = "Casey Flowers"
$nameMap["Thomas Spencer"] = "Wanda Price"
$nameMap["Garrett Williams"] = "Carmen Bell"
$nameMap["Jill Jackson"] = "Kimberly Suarez"
$nameMap["Bradley Sullivan"] = "Sonia Moore"
$nameMap["Jessica Moore"] = "Meghan Carpenter"
$nameMap["Mary Hernandez"] = "Audrey Berg"
$nameMap["Christopher Castaneda"] = "Leah Delacruz"
$nameMap["Danielle Wilkins"] = "Andrew Clark"
$nameMap["Cindy Pierce"] = "Jennifer Todd"
$nameMap["Steve Cooper"] = "Heidi Alvarez"
$nameMap["Jessica Bradshaw"] = "Kathleen Valenzuela"
$nameMap["Christian Richardson"] = "Jonathan Levine"
$nameMap["Dr. Norma Ramirez MD"] = "Amanda Weaver"
$nameMap["Dominique Valdez"] = "Dylan Foster"
$nameMap["Mallory Logan"] = "Guy Butler"
$nameMap["Robert Mills"] = "Christina Lara"
$nameMap["Gina Rios"] = "Whitney Hernandez"
$nameMap["Crystal Lopez"] = "Antonio Hale"
$nameMap["Deanna Cabrera"] = "Jonathon Caldwell"
$nameMap["Amanda Campbell"] = "Katrina Thompson"
$nameMap["Mary Reyes"] = "Thomas Martinez"
$nameMap["Matthew Doyle"] = "Melissa Berger"
$nameMap["Julie Ewing"] = "Nicholas Jennings"
$nameMap["Deborah Lang"] = "Carolyn Mejia"
$nameMap["Brian Spears"] = "Jeffrey Williams"
$nameMap["Ryan Bennett"] = "Jennifer Hayes"
$nameMap["Valerie Sanders"] = "Michael Paul"
$nameMap["Justin Spencer"] = "Doris Sims"
$nameMap["Lisa Mcbride"] = "Kristen Mcconnell"
$nameMap["Kenneth Owens"] = "Krystal Brooks"
$nameMap["Victoria Malone"] = "Nicholas White"
$nameMap["Aaron Hawkins"] = "Michelle Andersen"
$nameMap["Charles Harris"] = "Charlene Ramirez"
$nameMap["Kyle Conway"] = "Emily Lloyd"
$nameMap["Michelle Davis"] = "Steven Martin"
$nameMap["Sean Russell"] = "Jack Weber"
$nameMap["Leslie Callahan"] = "Thomas Johnson"
$nameMap["Monica Wallace"] = "Hector Rodriguez"
$nameMap["Bruce English"] = "Justin Shaw"
$nameMap["April Dawson"] = "Kevin Hernandez"
$nameMap["Kimberly Berger"] = "Jean Benitez"
$nameMap["George Banks"] = "Catherine Lopez"
$nameMap["Sheila Mendoza"] = "Jennifer Guzman"
$nameMap["Isabella Johnson"] = "Amanda Lewis"
$nameMap["Kathleen Gonzales"] = "Larry Williams"
$nameMap["Gabriela Jackson"] = "Joe Colon"
$nameMap["Victoria Frederick"] = "John Jensen"
$nameMap["Krystal Kerr"] = "John Willis"
$nameMap["Angela Velez"] = "Kimberly Gibbs"
$nameMap["Edward Conway"] = "Melissa Odonnell"
$nameMap["Kristine Smith"] = "Nathan Waters"
$nameMap["Katherine Cole"] = "Kimberly Fernandez"
$nameMap["David Wang"] = "Ronald Barr"
$nameMap["Paul Walter"] = "Christopher Kaiser"
$nameMap["Brenda Thompson"] = "Todd Price"
$nameMap["Ethan Tucker"] = "Jennifer Kim"
$nameMap["David Jackson"] = "Christopher Anderson"
$nameMap["Denise Carlson"] = "Kristi Clarke"
$nameMap["Christopher Flores Jr."] = "Jessica Cobb"
$nameMap["Paula Hanson"] = "Barbara Freeman"

$teamIdMap = @{}
$teamIdMap["9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"] = "1921564e-062f-4597-95dc-95074a120bb7"
$teamIdMap["61f0969e-22a4-4374-8588-d6511915b05e"] = "fb8d99d6-9f17-469f-9dd5-003d5c7968d9"
$teamIdMap["e7bb31c1-e095-453b-95ff-565ea62efb0a"] = "166314d0-c6c7-44e5-b2c2-2c69a7f25cd8"
$teamIdMap["0eeb011c-24fb-4476-91f7-d8e28ae49c2f"] = "0874d67c-598f-44e3-920b-4561908dfe3b"
$teamIdMap["e0228b4f-7807-45db-a3f6-8c6e1f4adf41"] = "2111d74d-0d3e-4e4f-9f1b-b6e8c0cfe0f4"
$teamIdMap["7fc75193-58a0-4e7d-ab42-382ec10a8be4"] = "9e9e4955-a3b5-41f2-a00c-bd1f5310c968"
$teamIdMap["5e1a20f6-82bf-4dee-aa79-41702d9feb41"] = "9e34f94d-c31f-4c23-99d2-a2822dd88d51"
$teamIdMap["7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"] = "32128616-e39a-4800-9f53-f18630bc75d0"
$teamIdMap["e5c40d19-b03a-4f5a-82c8-25540cd45e07"] = "a75e188a-12c5-4285-bcec-2d97ca3c7153"
$teamIdMap["0435a227-38e8-494e-b1bf-271b00893eae"] = "bf18fea4-73ce-4d76-86f9-bbf70bcbff6a"
$teamIdMap["630f61e8-543f-46e2-af63-2b62e8bc4fd2"] = "50be74bd-92cf-49cc-97dc-5599f0670e91"
$teamIdMap["3b8adc57-0f6f-482c-8306-9830e819d666"] = "1c9e5a86-8cbe-46d6-989c-5e819a6aaf27"
$teamIdMap["94742748-e7ab-454b-8ff6-9893440bd059"] = "d55d68ed-1702-4f23-b775-3777dd67494c"
$teamIdMap["5184566d-523a-4432-848d-ac234ffb6ac6"] = "2412507d-742b-4053-9f4e-e96905c184cd"
$teamIdMap["4279fd55-c2c1-440d-abaa-430f3c27be44"] = "92c24ef1-b25c-4ff5-bb92-9058df7f7b6e"
$teamIdMap["e5ca6e2b-5f54-4acd-ad7b-03e631313986"] = "a0768df2-6a70-45cb-9b17-170d678f0408"
$teamIdMap["6afc31f0-3916-443a-92c4-b5eb425a9bc3"] = "003c4550-974d-4aa8-9276-e3ef21048e79"
$teamIdMap["850a92da-c3d6-4fb9-a510-99626e9ad312"] = "19995c9e-d80e-471a-9813-a42cd6f71119"
$teamIdMap["57ebee16-96d2-46a6-ab16-2476b305fd91"] = "5c8943e1-fcd8-4718-b614-972b56a08435"
$teamIdMap["e42288a3-b5af-4464-bc45-85d438bcea11"] = "5c8fd640-4c38-43fd-be97-b526717d2df2"
$teamIdMap["de503c24-f17d-47a9-9a47-6f0a194f8c9c"] = "4c6a902c-c4c5-4ae2-88a1-05cd3ec101b8"
$teamIdMap["1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"] = "331efece-b4e3-4565-86c3-931d3f793263"
$teamIdMap["718c6b8f-7c00-4bcb-b53c-8f3f42154362"] = "af854dff-40e7-4c50-a0f1-6dd6958c6784"

$lastRow = $ws.UsedRange.Rows.Count
$updatedNames = 0
$updatedTeamIds = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $qCell = $ws.Cells.Item($r, 17)
    $qVal = $qCell.Value2
    if ($qVal -ne $null -and $nameMap.ContainsKey($qVal)) {
        $qCell.Value = $nameMap[$qVal]
        $updatedNames++
    }
}
for ($r = 2; $r -le $lastRow; $r++) {
    $sCell = $ws.Cells.Item($r, 19)
    $sVal = $sCell.Value2
    if ($sVal -ne $null -and $teamIdMap.ContainsKey($sVal)) {
        $sCell.Value = $teamIdMap[$sVal]
        $updatedTeamIds++
    }
}
Write-Host "Updated name cells:" $updatedNames
Write-Host "Updated teamid cells:" $updatedTeamIds
